$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column E ("VALOR") -------------------------------------------
# Header
$ws.Cells.Item(1, 5).Value = "VALOR"

# Data rows 2..110 get sequential values 1..109 (row number - 1)
$lastRow = 110
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $r - 1
}

# --- Highlight the rows for segments included in the primary motor area ---
# (ctx-lh-paracentral, ctx-lh-precentral, ctx-rh-paracentral, ctx-rh-precentral)
$highlightRows = @(58, 65, 92, 99)
foreach ($r in $highlightRows) {
    $ws.Range("A$r`:D$r").Interior.Color = 65535
    $ws.Range("E$r").Interior.Color = 65535
}

# --- View state -------------------------------------------------------------
$ws.Range("F1:F5").Select()
